$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (before former row 366),
# pushing the existing rows (366-432) down to (368-434) and growing the
# used range from A1:R432 to A1:R434.
$ws.Rows.Item(366).Insert()
$ws.Rows.Item(366).Insert()

# Populate the new row 366 (Primera) with this week's data.
$ws.Range("A366").Value = 3
$ws.Range("B366").Value = "Femacal de La Calera"
$ws.Range("C366").Value = "Coquimbo"
$ws.Range("D366").Value = 44504
$ws.Range("E366").Value = 5
$ws.Range("F366").Value = 100114014
$ws.Range("G366").Value = "Betarraga"
$ws.Range("H366").Value = "Sin especificar"
$ws.Range("I366").Value = "Primera"
$ws.Range("J366").Value = 3300
$ws.Range("K366").Value = 500
$ws.Range("L366").Value = 550
$ws.Range("M366").Value = 526
$ws.Range("N366").Value = "$/paquete 4 unidades"
$ws.Range("O366").Value = "Provincia de Quillota"
$ws.Range("P366").Value = 132
$ws.Range("Q366").Value = 4
$ws.Range("R366").Value = "Hortaliza"

# Populate the new row 367 (Segunda) with this week's data.
$ws.Range("A367").Value = 3
$ws.Range("B367").Value = "Femacal de La Calera"
$ws.Range("C367").Value = "Coquimbo"
$ws.Range("D367").Value = 44504
$ws.Range("E367").Value = 5
$ws.Range("F367").Value = 100114014
$ws.Range("G367").Value = "Betarraga"
$ws.Range("H367").Value = "Sin especificar"
$ws.Range("I367").Value = "Segunda"
$ws.Range("J367").Value = 1500
$ws.Range("K367").Value = 400
$ws.Range("L367").Value = 400
$ws.Range("M367").Value = 400
$ws.Range("N367").Value = "$/paquete 4 unidades"
$ws.Range("O367").Value = "Provincia de Quillota"
$ws.Range("P367").Value = 100
$ws.Range("Q367").Value = 4
$ws.Range("R367").Value = "Hortaliza"
